$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a "last changed" date serial for every data
# row (rows 2-369). The whole column of values moves forward by one day:
# 45171 (2023-09-02) -> 45172 (2023-09-03).
for ($r = 2; $r -le 369; $r++) {
    $ws.Cells.Item($r, 3).Value = 45172
}
